# Applies the "ggplot -> ggplot2" / "Music Taste -> Breakfast Habits" title
# edits described by the commit diff.
#
# Three text boxes change, identified by their (stable) shape names so the
# lookup does not depend on positional shape ordering:
#   Slide 12 - "Google Shape;183;p24"  -> "ggplot2 syntax"
#   Slide 13 - "Google Shape;192;p25"  -> "ggplot2 syntax"
#   Slide 37 - "Google Shape;362;p49"  -> "Task 5 <en dash> Does Age interact with Breakfast Habits?"

$p = $ppt.ActivePresentation

$enDash = [char]0x2013

function Set-ShapeText {
    param(
        [int]$SlideIndex,
        [string]$ShapeName,
        [string]$NewText
    )

    $slide = $p.Slides.Item($SlideIndex)
    $shape = $slide.Shapes.Item($ShapeName)
    $shape.TextFrame.TextRange.Text = $NewText
}

Set-ShapeText 12 "Google Shape;183;p24" "ggplot2 syntax"
Set-ShapeText 13 "Google Shape;192;p25" "ggplot2 syntax"
Set-ShapeText 37 "Google Shape;362;p49" ("Task 5 " + $enDash + " Does Age interact with Breakfast Habits?")
